# Update cryptocurrency price (D) and volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "37.690.97"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.66%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "2.072.54"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -2.00%  "
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "232.91"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.39%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "58.47"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +1.33%  "
$ws.Cells.Item(8, 5).Value = "  -0.02%  "
$ws.Cells.Item(9, 5).Value = "  +0.44%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.0781"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.20%  "
$ws.Cells.Item(11, 5).Value = "  +3.28%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "2.378.86"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -2.04%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "14.71"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.74%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "20.91"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.49%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.773"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -1.16%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "5.35"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.59%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "2.076.18"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.83%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "37.599.21"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.69%  "
$ws.Cells.Item(19, 5).Value = "  -1.80%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "71.27"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.22%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0833"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.08%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "228.00"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.21%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.12%  "
$ws.Cells.Item(24, 5).Value = "  -0.33%  "
$ws.Cells.Item(25, 5).Value = "  -2.70%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "171.18"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.08%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "9.02"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.75%  "
$ws.Cells.Item(28, 5).Value = "  -1.13%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "19.43"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.65%  "
$ws.Cells.Item(30, 5).Value = "  -2.55%  "
$ws.Cells.Item(31, 5).Value = "  +1.95%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "4.66"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.94%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.0632"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.31%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "4.65"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.26%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "2.48"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -5.34%  "
$ws.Cells.Item(36, 5).Value = "  -0.72%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "3.37"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.21%  "
$ws.Cells.Item(38, 5).Value = "  -0.02%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "5.32"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.61%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "99.69"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.85%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.0970"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.94%  "
$ws.Cells.Item(42, 5).Value = "  -2.01%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.0214"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.17%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "16.59"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +6.70%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "1.436.32"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.83%  "
$ws.Cells.Item(46, 5).Value = "  -0.51%  "
$ws.Cells.Item(47, 5).Value = "  +2.07%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "7.40"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +1.88%  "
$ws.Cells.Item(50, 5).Value = "  -1.70%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "2.263.54"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -2.10%  "
